$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -is [string]) {
            $nv = $v
            $nv = $nv.Replace("D64", "D69")
            $nv = $nv.Replace("D80", "D86")
            $nv = $nv.Replace("D51", "D55")
            $nv = $nv.Replace("S30", "S31")
            if ($nv -ne $v) {
                $cell.Value = $nv
            }
        }
    }
}
